$wb = $excel.ActiveWorkbook

# Map: sheet name -> (ResultProd A2, DateProd B2, DateDemo D2)
$updates = @{
    "CreateUser"    = @("Pass", "Fri Aug 22 22:33:22 IST 2025", "Wed Aug 20 22:56:56 IST 2025")
    "FindUser"      = @("Pass", "Fri Aug 22 22:33:57 IST 2025", "Wed Aug 20 22:57:31 IST 2025")
    "ModifyUser"    = @("Pass", "Fri Aug 22 22:35:06 IST 2025", "Thu Aug 21 00:23:23 IST 2025")
    "ModifyUserPwd" = @("Pass", "Mon Aug 25 20:03:36 IST 2025", "Mon Aug 25 19:53:00 IST 2025")
    "FindCaseUser"  = @("Pass", "Mon Aug 25 20:05:26 IST 2025", "Wed Aug 20 22:59:25 IST 2025")
    "AddDeleteRole" = @("Pass", "Fri Aug 22 22:32:10 IST 2025", "Wed Aug 20 22:55:55 IST 2025")
    "SearchRole"    = @("Pass", "Fri Aug 22 22:32:48 IST 2025", "Wed Aug 20 22:56:28 IST 2025")
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $updates[$sheetName]

    $ws.Range("A2").Value = $values[0]
    $ws.Range("B2").Value = $values[1]
    $ws.Range("D2").Value = $values[2]
}
